$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("27").Insert()
$ws.Rows("21").Copy()
$ws.Rows("27").PasteSpecial(-4122)
Write-Host "done"
